$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range('A2').Value = '航天动力'
$ws.Range('B2').Value = '摩尔线程-U'
$ws.Range('A3').Value = '雪人集团'
$ws.Range('C3').Value = '航天动力'
$ws.Range('A4').Value = '航天发展'
$ws.Range('B4').Value = '航天动力'
$ws.Range('C4').Value = '龙洲股份'
$ws.Range('A5').Value = '摩尔线程-U'
$ws.Range('B5').Value = '雪人集团'
$ws.Range('C5').Value = '东百集团'
$ws.Range('A6').Value = '永辉超市'
$ws.Range('B6').Value = '航天发展'
$ws.Range('C6').Value = '航天发展'
$ws.Range('A7').Value = '顺灏股份'
$ws.Range('B7').Value = '东百集团'
$ws.Range('C7').Value = '平潭发展'
$ws.Range('A8').Value = '华菱线缆'
$ws.Range('B8').Value = '中兴通讯'
$ws.Range('C8').Value = '雪人集团'
$ws.Range('B9').Value = '永鼎股份'
$ws.Range('C9').Value = '摩尔线程'
$ws.Range('A10').Value = '中超控股'
$ws.Range('B10').Value = '平潭发展'
$ws.Range('C10').Value = '顺灏股份'
$ws.Range('B11').Value = '龙洲股份'
$ws.Range('C11').Value = '博纳影业'
$ws.Range('A12').Value = '龙洲股份'
$ws.Range('B12').Value = '顺灏股份'
$ws.Range('C12').Value = '西部材料'
$ws.Range('A13').Value = '航天机电'
$ws.Range('B13').Value = '华菱线缆'
$ws.Range('C13').Value = '安妮股份'
$ws.Range('A14').Value = '永鼎股份'
$ws.Range('B14').Value = '中国西电'
$ws.Range('C14').Value = '再升科技'
$ws.Range('A15').Value = '中兴通讯'
$ws.Range('B15').Value = '安泰科技'
$ws.Range('C15').Value = '永鼎股份'
$ws.Range('A16').Value = '再升科技'
$ws.Range('B16').Value = '中能电气'
$ws.Range('C16').Value = '中超控股'
$ws.Range('A17').Value = '国机重装'
$ws.Range('B17').Value = '航天机电'
$ws.Range('C17').Value = '赛微电子'
$ws.Range('A18').Value = '安妮股份'
$ws.Range('B18').Value = '安妮股份'
$ws.Range('C18').Value = '航天机电'
$ws.Range('A19').Value = '通光线缆'
$ws.Range('B19').Value = '特变电工'
$ws.Range('C19').Value = '国机重装'
$ws.Range('A20').Value = '中能电气'
$ws.Range('B20').Value = '赛微电子'
$ws.Range('C20').Value = '中国西电'
$ws.Range('A21').Value = '中国西电'
$ws.Range('B21').Value = '国机重装'
$ws.Range('C21').Value = '雷科防务'
